# docs/resources/SourceDictionary.xlsx template update:
# "collection event.resource" + "collection event.name" columns are merged
# into a single "collection event" column on both the "Variables" sheet and
# the "Repeated variables" sheet.

$wb = $excel.ActiveWorkbook

# --- Variables sheet ---------------------------------------------------
# Old layout: ... H=references, I=collection event.resource,
#             J=collection event.name, K=description, ...
# New layout: ... H=references, I=collection event, J=description, ...
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Columns.Item(10).Delete()
$wsVariables.Range("I1").Value = "collection event"

# --- Repeated variables sheet -------------------------------------------
# Old layout: ... F=is repeat of.name, G=collection event.resource,
#             H=collection event.name, I=since version, J=until version
# New layout: ... F=is repeat of.name, G=collection event,
#             H=since version, I=until version
$wsRepeated = $wb.Worksheets.Item("Repeated variables")
$wsRepeated.Columns.Item(8).Delete()
$wsRepeated.Range("G1").Value = "collection event"

# --- selection / active cell bookkeeping --------------------------------
$wsRepeated.Activate() | Out-Null
$wsRepeated.Range("A2").Select() | Out-Null

$wsVariables.Activate() | Out-Null
$wsVariables.Range("A2").Select() | Out-Null

$wsDatasets = $wb.Worksheets.Item("Datasets")
$wsDatasets.Activate() | Out-Null
$wsDatasets.Range("A2").Select() | Out-Null
